$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8245447874069214
$ws.Range("B1").Value = 1.518949627876282
$ws.Range("C1").Value = 6.416202545166016
$ws.Range("D1").Value = 1.828590273857117
$ws.Range("E1").Value = 1.068598508834839
